# Update the statistical description table (Sheet1) with refreshed
# summary statistics (Mean, STD, quartiles, etc.) reflecting more data
# used in the underlying analysis.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

    $ws.Range("C2").Value = 555.3976844498427
    $ws.Range("D2").Value = 136.1086402776294
    $ws.Range("G2").Value = 516
    $ws.Range("H2").Value = 624
    $ws.Range("C3").Value = 40.65229553045877
    $ws.Range("D3").Value = 4.778124445764926
    $ws.Range("F3").Value = 37.68
    $ws.Range("H3").Value = 43.32
    $ws.Range("C4").Value = 1.392464612725405
    $ws.Range("D4").Value = 2.248701706053477
    $ws.Range("H4").Value = 1.8
    $ws.Range("C5").Value = 323.9328607341399
    $ws.Range("D5").Value = 10.2513172701727
    $ws.Range("F5").Value = 318.07
    $ws.Range("G5").Value = 325.76
    $ws.Range("H5").Value = 332.35
    $ws.Range("C6").Value = 21.00981532631038
    $ws.Range("D6").Value = 2.124086112344261
    $ws.Range("E6").Value = 14.76
    $ws.Range("F6").Value = 19.46
    $ws.Range("G6").Value = 20.64
    $ws.Range("C7").Value = -76.33482824130306
    $ws.Range("D7").Value = 22.44113227402011
    $ws.Range("C8").Value = 7.762797186990312
    $ws.Range("D8").Value = 6.838123600698207
    $ws.Range("C9").Value = 9.322399872309944
    $ws.Range("D9").Value = 1.688394112346986
    $ws.Range("C10").Value = 867.8303629823408
    $ws.Range("D10").Value = 0.4610342240256628
    $ws.Range("C11").Value = 0.5569113180728144
    $ws.Range("D11").Value = 0.5906378012529424
    $ws.Range("C12").Value = 22.6893627954779
    $ws.Range("D12").Value = 12.27781762319351
    $ws.Range("C13").Value = 0.6716720779220779
    $ws.Range("D13").Value = 0.7482962772579064
    $ws.Range("C14").Value = 1.825914268585132
    $ws.Range("D14").Value = 1.665830696243136
    $ws.Range("C15").Value = 93.73482824130286
    $ws.Range("D15").Value = 22.4411322740201
    $ws.Range("C16").Value = -85.65015002324094
    $ws.Range("D16").Value = 20.23861619242166
    $ws.Range("F16").Value = -101.1469917995764
    $ws.Range("H16").Value = -67.68978441047734
    $ws.Range("C17").Value = -77.88735283625061
    $ws.Range("D17").Value = 24.81704205386455
    $ws.Range("F17").Value = -92.2376019773414
    $ws.Range("G17").Value = -75.41392685158225
    $ws.Range("H17").Value = -56.79009749652566
